# Apply updated dSF (column F) values after repulling data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    9  = 1
    13 = 0
    17 = 0
    23 = -3
    25 = -2
    28 = -3
    30 = -3
    31 = -4
    36 = -6
    39 = -2
    44 = -1
    47 = 0
    51 = -4
    53 = -1
    55 = -2
    63 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
